$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the header typo: "Thickness of the head (Max))" -> "Thickness of the head (Max)"
$ws.Cells.Item(1, 15).Value = "Thickness of the head (Max)"

# Update the Standards column (B) values from "ISO-4014-2011" to "ISO 4014" for all data rows
for ($r = 2; $r -le 39; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value2 -eq "ISO-4014-2011") {
        $cell.Value = "ISO 4014"
    }
}

# Update the sheet view: remove frozen/scrolled topLeftCell, change selection to F8:F9 (active F9)
$ws.Range("A1").Select() | Out-Null
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("F8:F9").Select() | Out-Null
